$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (thin box border) for the new rows 40-51 by reusing the existing
# bordered style from row 39 (same style as all other data rows).
$ws.Range("A39:B39").Copy()
$ws.Range("A40:B51").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Write the full, reordered/updated Variable/Description table (rows 1-51).
$ws.Range("A1").Value = 'Variable'
$ws.Range("B1").Value = 'Description'
$ws.Range("A2").Value = 'alias_sequencing_id'
$ws.Range("B2").Value = 'De-identified sample identification number'
$ws.Range("A3").Value = 'alias_study_id'
$ws.Range("B3").Value = 'De-identified subject identification number'
$ws.Range("A4").Value = 'SampleTiming'
$ws.Range("B4").Value = 'Timing of sample collection relative to SARS-CoV-2 infection or exposure (acute for all samples in this analysis)'
$ws.Range("A5").Value = 'SampleType'
$ws.Range("B5").Value = 'np = nasopharyngeal, pax = blood, nasal = nasal'
$ws.Range("A6").Value = 'year'
$ws.Range("B6").Value = 'Year of sample collection'
$ws.Range("A7").Value = 'study'
$ws.Range("B7").Value = 'Study into which subject was enrolled (brave or messi)'
$ws.Range("A8").Value = 'hospital'
$ws.Range("B8").Value = 'Did the subject require hospitalization? (Y/N)'
$ws.Range("A9").Value = 'batch_num'
$ws.Range("B9").Value = 'Sequencing batch number'
$ws.Range("A10").Value = 'RIN'
$ws.Range("B10").Value = 'Sample RNA integrity number'
$ws.Range("A11").Value = 'age'
$ws.Range("B11").Value = 'Age in years'
$ws.Range("A12").Value = 'age_cat'
$ws.Range("B12").Value = 'Age category (0-5 years, 6-13 years, 14-20 years, Adult)'
$ws.Range("A13").Value = 'sex'
$ws.Range("B13").Value = 'F = female, M = male'
$ws.Range("A14").Value = 'race'
$ws.Range("B14").Value = 'Race and ethncity data'
$ws.Range("A15").Value = 'hispanic'
$ws.Range("B15").Value = 'Hispanic ethnicity (Yes/No)'
$ws.Range("A16").Value = 'corona'
$ws.Range("B16").Value = 'SARS-CoV-2 infection status (Negative, Positive)'
$ws.Range("A17").Value = 'symptoms'
$ws.Range("B17").Value = 'Presence of any symptoms (Y/N)'
$ws.Range("A18").Value = 'group'
$ws.Range("B18").Value = 'SARS-CoV-2 infection status and symptoms (NEG_ASX = asymptomatic uninfected, POS_ASX = asymptomatic infected, POS_SX = symptomatic infected)'
$ws.Range("A19").Value = 'clinical_pcr'
$ws.Range("B19").Value = 'Results of PCR performed for clinical care'
$ws.Range("A20").Value = 'research_pcr'
$ws.Range("B20").Value = 'Results of quantitative PCR performed for research'
$ws.Range("A21").Value = 'vl_copies'
$ws.Range("B21").Value = 'Viral load (copies/mL)'
$ws.Range("A22").Value = 'vaccine_doses'
$ws.Range("B22").Value = 'Number of COVID-19 vaccine doses received prior to sample collection'
$ws.Range("A23").Value = 'timing_sx'
$ws.Range("B23").Value = 'Timing of sample relative to symptom onset in days (NA unless symptomatic infected)'
$ws.Range("A24").Value = 'timing_dx'
$ws.Range("B24").Value = 'Timing of sample relative to SARS-CoV-2 diagnosis in days (NA unless infected)'
$ws.Range("A25").Value = 'obesity'
$ws.Range("B25").Value = 'BMI >95% for age (Y/N); missing for children <2 years of age'
$ws.Range("A26").Value = 'comorbidity_oth'
$ws.Range("B26").Value = 'Presence of a comorbid medical condition other than obesity (Y/N)'
$ws.Range("A27").Value = 'asthma'
$ws.Range("B27").Value = 'History of physician-diagnosed asthma (Y/N)'
$ws.Range("A28").Value = 'pulm_oth'
$ws.Range("B28").Value = 'History of other pulmonary condition (Y/N)'
$ws.Range("A29").Value = 'htn'
$ws.Range("B29").Value = 'History of hypertension (Y/N)'
$ws.Range("A30").Value = 'cardiac_oth'
$ws.Range("B30").Value = 'History of other cardiac condition (Y/N)'
$ws.Range("A31").Value = 'diabetes'
$ws.Range("B31").Value = 'History of diabetes mellitus (Y/N)'
$ws.Range("A32").Value = 'neuro'
$ws.Range("B32").Value = 'History of a chronic neurological disorder (Y/N)'
$ws.Range("A33").Value = 'renal'
$ws.Range("B33").Value = 'History of a chronic renal disorder (Y/N)'
$ws.Range("A34").Value = 'cancer'
$ws.Range("B34").Value = 'History of malignancy (Y/N)'
$ws.Range("A35").Value = 'immuno'
$ws.Range("B35").Value = 'Immunosuppressed status (Y/N)'
$ws.Range("A36").Value = 'fever'
$ws.Range("B36").Value = 'Presence of fever (Y/N)'
$ws.Range("A37").Value = 'cough'
$ws.Range("B37").Value = 'Presence of cough (Y/N)'
$ws.Range("A38").Value = 'sob'
$ws.Range("B38").Value = 'Presence of shortness of breath (Y/N)'
$ws.Range("A39").Value = 'sorethroat'
$ws.Range("B39").Value = 'Presence of sore throat (Y/N)'
$ws.Range("A40").Value = 'rhinorrhea'
$ws.Range("B40").Value = 'Presence of rhinorrhea (Y/N)'
$ws.Range("A41").Value = 'congestion'
$ws.Range("B41").Value = 'Presence of congestion (Y/N)'
$ws.Range("A42").Value = 'headache'
$ws.Range("B42").Value = 'Presence of headache (Y/N)'
$ws.Range("A43").Value = 'abd_pain'
$ws.Range("B43").Value = 'Presence of abdominal pain (Y/N)'
$ws.Range("A44").Value = 'diarrhea'
$ws.Range("B44").Value = 'Presence of diarrhea (Y/N)'
$ws.Range("A45").Value = 'anosmia'
$ws.Range("B45").Value = 'Loss of smell (Y/N)'
$ws.Range("A46").Value = 'dysgeusia'
$ws.Range("B46").Value = 'Loss of taste (Y/N)'
$ws.Range("A47").Value = 'chestpain'
$ws.Range("B47").Value = 'Presence of chest pain (Y/N)'
$ws.Range("A48").Value = 'myalgias'
$ws.Range("B48").Value = 'Presence of myalgias (Y/N)'
$ws.Range("A49").Value = 'joint_pain'
$ws.Range("B49").Value = 'Presence of joint pains (Y/N)'
$ws.Range("A50").Value = 'sars2_variant'
$ws.Range("B50").Value = 'SARS-CoV-2 lineage based on genomic seqencing'
$ws.Range("A51").Value = 'neut_ID50_2mo'
$ws.Range("B51").Value = 'Neutralizing activity of serum antibodies to SARS-CoV-2 (D614G) at 2 months after acute infection'

# Restore the selection to match the saved workbook state (active cell B8).
$ws.Range("B8").Select()

Write-Output "BRAVE RNASeq data dictionary updated"
